# feat: add 2022-Q1 data
#
# 1) Insert a new sheet "2022-Q1" (between "2021-Q1" and "总计") holding the
#    per-fund holdings for 2022-Q1.
# 2) Insert a new top row into "总计" summarising the 2022-Q1 quarter, above
#    the existing "2021-Q1" summary row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: build the new "2022-Q1" sheet by duplicating "总计" (so it inherits
# the same header/first-column cell style) immediately after "2021-Q1", then
# overwrite its contents.
#
# NOTE: worksheet variables resolve by *position*, not identity, so any
# handle obtained before a sheet add/copy/delete goes stale afterwards (it
# silently starts pointing at whatever sheet now sits at that index). Always
# re-fetch `Worksheets.Item(name)` fresh right after such a structural
# change instead of reusing an older reference.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("总计").Copy($null, $wb.Worksheets.Item("2021-Q1"))
# The copy lands immediately after "2021-Q1", i.e. at index 2 (whatever
# auto-generated name it got, e.g. "总计 (2)") - grab it positionally rather
# than relying on the exact auto-generated name.
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"

# Extend the header row from B:D to B:H, copying the existing header style
# across (same-sheet paste-special preserves formatting in this engine).
$newSheet.Range("B1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)

# Add two more data rows (the sheet starts with just one data row, copied
# from "总计"); each new A-cell picks up the existing A2 style.
$newSheet.Rows.Item(3).Insert()
$newSheet.Rows.Item(3).Insert()
$newSheet.Range("A2").Copy()
$newSheet.Range("A3:A4").PasteSpecial(-4122)

# Header text.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row index column (0,1,2).
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2

# Fund-holding rows (B:G are kept as text, matching the source data which
# preserves trailing zeros such as "350.10"; H is numeric).
$newSheet.Range("B2:G4").NumberFormat = "@"

$newSheet.Range("B2").Value = "513050"
$newSheet.Range("C2").Value = "易方达中证海外中国互联网50 QDII-ETF"
$newSheet.Range("D2").Value = "350.10"
$newSheet.Range("E2").Value = "98.05"
$newSheet.Range("F2").Value = "4.67"
$newSheet.Range("G2").Value = "16.3497"
$newSheet.Range("H2").Value = 5

$newSheet.Range("B3").Value = "159605"
$newSheet.Range("C3").Value = "广发中证海外中国互联网30（QDII-ETF）"
$newSheet.Range("D3").Value = "29.04"
$newSheet.Range("E3").Value = "98.61"
$newSheet.Range("F3").Value = "7.09"
$newSheet.Range("G3").Value = "2.0589"
$newSheet.Range("H3").Value = 5

$newSheet.Range("B4").Value = "159607"
$newSheet.Range("C4").Value = "嘉实中证海外中国互联网30ETF（QDII）"
$newSheet.Range("D4").Value = "5.79"
$newSheet.Range("E4").Value = "98.25"
$newSheet.Range("F4").Value = "7.14"
$newSheet.Range("G4").Value = "0.4134"
$newSheet.Range("H4").Value = 5

# Drop the "force text" number format again so the cells end up unstyled,
# same as the source workbook (only the A column + header row carry a style).
$newSheet.Range("B2:G4").ClearFormats()

# ---------------------------------------------------------------------------
# Step 2: insert the 2022-Q1 summary row at the top of "总计" (re-fetched by
# name since the sheet copy above shifted its position from 2 to 3).
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 18.82
$wsTotal.Range("B2:D2").ClearFormats()

# The pre-existing "2021-Q1" row (pushed from row 2 to row 3) is now the
# second entry, so its running index bumps from 0 to 1.
$wsTotal.Range("A3").Value = 1

# A2 should carry the same style as the rest of the A column (copied from
# A3, the old A2, which already has it).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

# Restore the original active sheet/selection (unaffected by this edit).
$wb.Worksheets.Item("2021-Q1").Activate()
[void]$wb.Worksheets.Item("2021-Q1").Range("A1").Select()
